# Modification of data preprocessing
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated accuracy numbers for the existing model rows (C3:C12) after
# re-running the preprocessing step.
$ws.Range("C3").Value  = 71.550196435691205
$ws.Range("C4").Value  = 94.864480171692904
$ws.Range("C5").Value  = 96.131891332527104
$ws.Range("C6").Value  = 96.136086308924206
$ws.Range("C7").Value  = 90.313239586760503
$ws.Range("C8").Value  = 69.781725315391697
$ws.Range("C9").Value  = 89.384114536514403
$ws.Range("C10").Value = 95.496711032342702
$ws.Range("C11").Value = 95.4664193060201
$ws.Range("C12").Value = 97.665839938404105

# New trailing row: label + average of the accuracy column.
$ws.Range("B13").Value = "Avg"

# Give B13 the same "left border only" look used by the model-name column
# (B4:B11), then flip its fill so the style is tracked as explicitly
# formatted (matches the new cellXfs entry added for this row).
$ws.Range("B4").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("B13").Interior.Pattern = 1

$ws.Range("C13").Formula = "=AVERAGE(C3:C12)"

# Match the saved selection/view state.
$ws.Range("E10:E11").Select()
